function Set-TextValue($range, $value) {
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "29.351.72"
Set-TextValue $ws.Range("E2") "  -0.05%  "
Set-TextValue $ws.Range("D3") "1.879.16"
Set-TextValue $ws.Range("E3") "  +0.23%  "
Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("E4") "  +0.07%  "
Set-TextValue $ws.Range("D5") "0.7113"
Set-TextValue $ws.Range("E5") "  -0.21%  "
Set-TextValue $ws.Range("D6") "242.69"
Set-TextValue $ws.Range("E6") "  +0.34%  "
Set-TextValue $ws.Range("D7") "1.001"
Set-TextValue $ws.Range("E7") "  +0.12%  "
Set-TextValue $ws.Range("D8") "0.08054"
Set-TextValue $ws.Range("E8") "  +3.26%  "
Set-TextValue $ws.Range("D9") "0.3182"
Set-TextValue $ws.Range("E9") "  +2.19%  "
Set-TextValue $ws.Range("E10") "  -0.47%  "
Set-TextValue $ws.Range("D11") "0.08332"
Set-TextValue $ws.Range("E11") "  -1.32%  "
Set-TextValue $ws.Range("D12") "1.891.75"
Set-TextValue $ws.Range("E12") "  +0.89%  "
Set-TextValue $ws.Range("D13") "5.264"
Set-TextValue $ws.Range("E13") "  +0.36%  "
Set-TextValue $ws.Range("D14") "94.89"
Set-TextValue $ws.Range("E14") "  +4.09%  "
Set-TextValue $ws.Range("D15") "0.7177"
Set-TextValue $ws.Range("E15") "  +0.50%  "
Set-TextValue $ws.Range("D16") "6.391"
Set-TextValue $ws.Range("E16") "  +5.13%  "
Set-TextValue $ws.Range("D17") "0.000008648"
Set-TextValue $ws.Range("E17") "  +4.84%  "
Set-TextValue $ws.Range("D18") "29.354.30"
Set-TextValue $ws.Range("E18") "  -0.05%  "
Set-TextValue $ws.Range("D19") "242.92"
Set-TextValue $ws.Range("E19") "  +0.78%  "
Set-TextValue $ws.Range("D20") "13.33"
Set-TextValue $ws.Range("E20") "  +0.53%  "
Set-TextValue $ws.Range("D21") "2.134.54"
Set-TextValue $ws.Range("E21") "  +0.43%  "
Set-TextValue $ws.Range("D22") "1.001"
Set-TextValue $ws.Range("E22") "  +0.10%  "
Set-TextValue $ws.Range("D23") "7.825"
Set-TextValue $ws.Range("E23") "  +0.49%  "
Set-TextValue $ws.Range("D24") "1.002"
Set-TextValue $ws.Range("E24") "  +0.14%  "
Set-TextValue $ws.Range("D25") "0.1574"
Set-TextValue $ws.Range("E25") "  -1.36%  "
Set-TextValue $ws.Range("D26") "9.093"
Set-TextValue $ws.Range("E26") "  +0.26%  "
Set-TextValue $ws.Range("D27") "163.17"
Set-TextValue $ws.Range("E27") "  +0.02%  "
Set-TextValue $ws.Range("D28") "18.60"
Set-TextValue $ws.Range("E28") "  +0.18%  "
Set-TextValue $ws.Range("E29") "  -0.23%  "
Set-TextValue $ws.Range("E30") "  +0.26%  "
Set-TextValue $ws.Range("D31") "4.330"
Set-TextValue $ws.Range("E31") "  +0.00%  "
Set-TextValue $ws.Range("D32") "1.200"
Set-TextValue $ws.Range("E32") "  -6.85%  "
Set-TextValue $ws.Range("D33") "0.05411"
Set-TextValue $ws.Range("E33") "  +2.02%  "
Set-TextValue $ws.Range("E34") "  +0.26%  "
Set-TextValue $ws.Range("D35") "0.7734"
Set-TextValue $ws.Range("E35") "  +3.85%  "
Set-TextValue $ws.Range("D36") "1.190"
Set-TextValue $ws.Range("E36") "  +0.86%  "
Set-TextValue $ws.Range("E37") "  -0.44%  "
Set-TextValue $ws.Range("D38") "0.01892"
Set-TextValue $ws.Range("E38") "  +1.16%  "
Set-TextValue $ws.Range("D39") "1.266.49"
Set-TextValue $ws.Range("E39") "  +3.05%  "
Set-TextValue $ws.Range("D40") "2.753"
Set-TextValue $ws.Range("E40") "  +0.78%  "
Set-TextValue $ws.Range("D41") "6.500"
Set-TextValue $ws.Range("E41") "  -0.34%  "
Set-TextValue $ws.Range("D42") "114.04"
Set-TextValue $ws.Range("E42") "  +2.58%  "
Set-TextValue $ws.Range("D43") "74.70"
Set-TextValue $ws.Range("E43") "  +2.32%  "
Set-TextValue $ws.Range("D44") "0.9075"
Set-TextValue $ws.Range("E44") "  +1.71%  "
Set-TextValue $ws.Range("D45") "0.00000000132"
Set-TextValue $ws.Range("E45") "  +7.21%  "
Set-TextValue $ws.Range("D46") "1.001"
Set-TextValue $ws.Range("E46") "  +0.12%  "
Set-TextValue $ws.Range("D47") "2.030.88"
Set-TextValue $ws.Range("E47") "  +0.53%  "
Set-TextValue $ws.Range("D48") "1.809"
Set-TextValue $ws.Range("E48") "  -0.29%  "
Set-TextValue $ws.Range("D49") "0.5221"
Set-TextValue $ws.Range("E49") "  +0.13%  "
Set-TextValue $ws.Range("D50") "9.525"
Set-TextValue $ws.Range("E50") "  +0.86%  "
Set-TextValue $ws.Range("D51") "0.4379"
Set-TextValue $ws.Range("E51") "  +1.26%  "
